$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume data (keeping values as plain text,
# matching the original worksheet formatting which stores these as text strings).
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "244.43"
Set-TextValue "E2" "-0.95%"
Set-TextValue "D3" "27.10"
Set-TextValue "E3" "3.37%"
Set-TextValue "D4" "5.137"
Set-TextValue "E4" "1.49%"
Set-TextValue "D5" "0.05657"
Set-TextValue "E5" "1.01%"
Set-TextValue "D6" "6.472"
Set-TextValue "D7" "0.8211"
Set-TextValue "E7" "0.99%"
Set-TextValue "D8" "0.8405"
Set-TextValue "E8" "-0.14%"
Set-TextValue "D9" "0.1329"
Set-TextValue "E9" "-1.00%"
Set-TextValue "D10" "0.06900"
Set-TextValue "E10" "-0.81%"
Set-TextValue "D11" "0.02978"
Set-TextValue "E11" "6.47%"
Set-TextValue "D12" "0.09393"
Set-TextValue "E12" "-0.05%"
Set-TextValue "D13" "0.001522"
Set-TextValue "E13" "0.81%"
Set-TextValue "D14" "0.04265"
Set-TextValue "E14" "-8.71%"
Set-TextValue "D15" "0.0005981"
Set-TextValue "E15" "0.23%"
Set-TextValue "D16" "0.006141"
Set-TextValue "E16" "-0.42%"
Set-TextValue "D17" "3.515"
Set-TextValue "E17" "-1.23%"
Set-TextValue "D18" "3.004"
Set-TextValue "E18" "-1.25%"
Set-TextValue "D19" "2.313"
Set-TextValue "E19" "9.18%"
Set-TextValue "E20" "-1.36%"
Set-TextValue "E21" "-0.21%"
Set-TextValue "E22" "-3.48%"
Set-TextValue "D23" "3.578"
Set-TextValue "E23" "-4.31%"
Set-TextValue "E24" "-0.04%"
Set-TextValue "E25" "-2.01%"
Set-TextValue "D26" "0.004466"
Set-TextValue "E26" "-3.12%"
Set-TextValue "D27" "0.00009801"
Set-TextValue "E27" "2.14%"
Set-TextValue "D28" "0.00007259"
Set-TextValue "E28" "-47.74%"
Set-TextValue "D40" "0.03651"
Set-TextValue "E40" "-0.22%"
Set-TextValue "E41" "-1.58%"
Set-TextValue "D42" "0.1053"
Set-TextValue "E42" "-0.21%"
Set-TextValue "D43" "0.002300"
Set-TextValue "E43" "-11.32%"
Set-TextValue "D44" "0.008976"
Set-TextValue "E44" "2.46%"
Set-TextValue "D45" "0.00005370"
Set-TextValue "E45" "1.50%"
Set-TextValue "E46" "0.05%"
Set-TextValue "D47" "0.1010"
Set-TextValue "E47" "-36.80%"
Set-TextValue "D48" "0.002654"
Set-TextValue "E48" "28.84%"
Set-TextValue "E49" "0.05%"
Set-TextValue "E50" "0.05%"
